$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text format so numeric-looking strings
# (e.g. "1.00", "3.115.67") are preserved exactly as text, matching the
# original inline-string cell contents instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '63.419.09'
$ws.Range('E2').Value = '  -2.64%  '
$ws.Range('D3').Value = '3.123.06'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '557.74'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = '139.80'
$ws.Range('E6').Value = '  -6.32%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '3.117.22'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('D11').Value = '0.161'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '0.461'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').Value = '35.52'
$ws.Range('E13').Value = '  -5.31%  '
$ws.Range('E14').Value = '  -3.28%  '
$ws.Range('D15').Value = '3.629.61'
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '63.428.82'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('D17').Value = '0.112'
$ws.Range('D18').Value = '3.115.97'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').Value = '511.39'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '6.76'
$ws.Range('E20').Value = '  -1.44%  '
$ws.Range('D21').Value = '13.68'
$ws.Range('E21').Value = '  -3.11%  '
$ws.Range('D22').Value = '0.714'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').Value = '7.33'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = '12.49'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '78.42'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').Value = '8.37'
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  -8.88%  '
$ws.Range('D31').Value = '26.56'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').Value = '2.55'
$ws.Range('E32').Value = '  -6.94%  '
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('D34').Value = '59.85'
$ws.Range('E34').Value = '  +12.56%  '
$ws.Range('D35').Value = '536.96'
$ws.Range('E35').Value = '  -10.62%  '
$ws.Range('D36').Value = '6.01'
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('D37').Value = '5.28'
$ws.Range('E37').Value = '  -7.19%  '
$ws.Range('D38').Value = '0.0417'
$ws.Range('E38').Value = '  -3.74%  '
$ws.Range('D39').Value = '0.0802'
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('D40').Value = '3.078.98'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').Value = '0.121'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('D42').Value = '2.76'
$ws.Range('E42').Value = '  -8.37%  '
$ws.Range('D43').Value = '8.17'
$ws.Range('E43').Value = '  -2.54%  '
$ws.Range('D44').Value = '0.258'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D46').Value = '2.10'
$ws.Range('E46').Value = '  -5.48%  '
$ws.Range('D47').Value = '122.82'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '24.49'
$ws.Range('E48').Value = '  -6.71%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.108'
$ws.Range('E49').Value = '  -2.05%  '
$ws.Range('B50').Value = 'PEPE'
$ws.Range('C50').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D50').Value = '0.0₃0516'
$ws.Range('E50').Value = '  -6.51%  '
$ws.Range('D51').Value = '2.45'
$ws.Range('E51').Value = '  +66.08%  '

# Restore default (Normal) style so no stray style index is left on the cells,
# keeping the workbook formatting identical to the original.
$ws.Range("D2:E51").Style = "Normal"
